$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.10887317818423
$ws.Range("C2").Value = 11.3650495583562
$ws.Range("E2").Value = 16.59571206078571
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 25.21245522896909
$ws.Range("H2").Value = 12.54707813218811

$ws.Range("B3").Value = 14.30894659641013
$ws.Range("C3").Value = 10.71428075818618
$ws.Range("E3").Value = 15.64325292521643
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 24.82128794009859
$ws.Range("H3").Value = 12.5995419000182

$ws.Range("B4").Value = 13.79642824505341
$ws.Range("C4").Value = 10.29214069113385
$ws.Range("E4").Value = 15.0336158413651
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 24.60066255999779
$ws.Range("H4").Value = 12.63716491894867

$ws.Range("B5").Value = 13.58243340822524
$ws.Range("C5").Value = 10.11449612797202
$ws.Range("E5").Value = 14.77920870365523
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 24.51579785755388
$ws.Range("H5").Value = 12.65384129831487

$ws.Range("B6").Value = 13.54659693162307
$ws.Range("C6").Value = 10.08466091329749
$ws.Range("E6").Value = 14.73661256982682
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 24.50201387455851
$ws.Range("H6").Value = 12.6566911480946

$ws.Range("B7").Value = 13.79356270792069
$ws.Range("C7").Value = 10.28976758190879
$ws.Range("E7").Value = 15.03020861819096
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 24.59949748180511
$ws.Range("H7").Value = 12.63738439904023

$ws.Range("B8").Value = 14.83761819200115
$ws.Range("C8").Value = 11.14537664815471
$ws.Range("E8").Value = 16.27260058222469
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 25.07360442045384
$ws.Range("H8").Value = 12.56403545185865

$ws.Range("B9").Value = 16.70730653922966
$ws.Range("C9").Value = 12.64242998345016
$ws.Range("E9").Value = 18.62411721303596
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 26.1519113437168
$ws.Range("H9").Value = 12.46380737966844

$ws.Range("B10").Value = 17.96392978987947
$ws.Range("C10").Value = 13.6305202902591
$ws.Range("E10").Value = 20.29484425364355
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 27.02472522888752
$ws.Range("H10").Value = 12.41765636949583

$ws.Range("B11").Value = 18.50888775218222
$ws.Range("C11").Value = 14.05558663438246
$ws.Range("E11").Value = 21.01289937957582
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 27.43691421910216
$ws.Range("H11").Value = 12.40281112229978

$ws.Range("B12").Value = 18.71133087111512
$ws.Range("C12").Value = 14.21302969089083
$ws.Range("E12").Value = 21.27881432123309
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 27.59496987122466
$ws.Range("H12").Value = 12.39808778519321

$ws.Range("B13").Value = 18.6679066378776
$ws.Range("C13").Value = 14.17927816683989
$ws.Range("E13").Value = 21.22181089566174
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 27.56084554593762
$ws.Range("H13").Value = 12.39906487574467

$ws.Range("B14").Value = 18.52562186033195
$ws.Range("C14").Value = 14.06861015477792
$ws.Range("E14").Value = 21.03489633335404
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 27.44987930088805
$ws.Range("H14").Value = 12.40240445089227

$ws.Range("B15").Value = 18.43795556896152
$ws.Range("C15").Value = 14.00036417525759
$ws.Range("E15").Value = 20.91962590789258
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 27.38215941528273
$ws.Range("H15").Value = 12.40456741723979

$ws.Range("B16").Value = 17.92777099701518
$ws.Range("C16").Value = 13.60224894186689
$ws.Range("E16").Value = 20.2470764183676
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 26.99807461247434
$ws.Range("H16").Value = 12.41875149003098

$ws.Range("B17").Value = 17.60788942251626
$ws.Range("C17").Value = 13.35175560736382
$ws.Range("E17").Value = 19.82376731799244
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 26.76618156984216
$ws.Range("H17").Value = 12.42903807504735

$ws.Range("B18").Value = 17.42139390634369
$ws.Range("C18").Value = 13.20538002183271
$ws.Range("E18").Value = 19.57633873380563
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 26.63424079744398
$ws.Range("H18").Value = 12.43553263269658

$ws.Range("B19").Value = 17.35782178669836
$ws.Range("C19").Value = 13.1554253491349
$ws.Range("E19").Value = 19.4918836854512
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 26.58982120684303
$ws.Range("H19").Value = 12.4378304094971

$ws.Range("B20").Value = 17.64220159841066
$ws.Range("C20").Value = 13.37865900889039
$ws.Range("E20").Value = 19.86923807651075
$ws.Range("F20").Value = 20.2495528364879
$ws.Range("G20").Value = 26.79071965984567
$ws.Range("H20").Value = 12.42788312744717

$ws.Range("B21").Value = 18.56752133139301
$ws.Range("C21").Value = 14.10121161876771
$ws.Range("E21").Value = 21.08996012578658
$ws.Range("F21").Value = 21.46857628470567
$ws.Range("G21").Value = 27.48242102190673
$ws.Range("H21").Value = 12.40139905082014

$ws.Range("B22").Value = 19.14939175136637
$ws.Range("C22").Value = 14.5529214648925
$ws.Range("E22").Value = 21.85283520027986
$ws.Range("F22").Value = 22.22866616901555
$ws.Range("G22").Value = 27.94586500857743
$ws.Range("H22").Value = 12.3893317544469

$ws.Range("B23").Value = 18.84095320961656
$ws.Range("C23").Value = 14.31371466002683
$ws.Range("E23").Value = 21.44885871735851
$ws.Range("F23").Value = 21.82633154475864
$ws.Range("G23").Value = 27.69754452525425
$ws.Range("H23").Value = 12.39528825563114

$ws.Range("B24").Value = 17.62669712864225
$ws.Range("C24").Value = 13.36650334253616
$ws.Range("E24").Value = 19.848693403368
$ws.Range("F24").Value = 20.22900810905294
$ws.Range("G24").Value = 26.779621682536
$ws.Range("H24").Value = 12.42840347169196

$ws.Range("B25").Value = 16.22154101703097
$ws.Range("C25").Value = 12.25699011043208
$ws.Range("E25").Value = 17.97139014462723
$ws.Range("F25").Value = 18.34778573295697
$ws.Range("G25").Value = 25.84527235897597
$ws.Range("H25").Value = 12.48615799805519
